$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.421.80'
$ws.Range('E2').Value = '  -0.08%  '
$ws.Range('D3').Value = '1.848.53'
$ws.Range('E3').Value = '  -0.10%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9989'
$ws.Range('D4').ClearFormats()
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.85'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.93%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6318'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -3.69%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '3.259.17'
$ws.Range('E8').Value = '  +76.27%  '
$ws.Range('E9').Value = '  +1.34%  '
$ws.Range('E10').Value = '  -0.99%  '
$ws.Range('E11').Value = '  +1.10%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07714'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +1.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.988'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.54%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6852'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.10%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000009991'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +4.89%  '
$ws.Range('E16').Value = '  -0.97%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.180'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.91%  '
$ws.Range('D18').Value = '29.464.96'
$ws.Range('E18').Value = '  +0.00%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '231.68'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.43%  '
$ws.Range('E20').Value = '  -0.38%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.569'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.25%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.000'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '155.21'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1388'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -2.54%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.439'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.51%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.66'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.81%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.473'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.75%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.05810'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -3.67%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.258'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.34%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.129'
$ws.Range('D31').ClearFormats()
$ws.Range('E32').Value = '  -1.24%  '
$ws.Range('D33').Value = '3.401.70'
$ws.Range('E33').Value = '  +69.98%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.871'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.79%  '
$ws.Range('E35').Value = '  -1.55%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7195'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.32%  '
$ws.Range('D38').Value = '1.248.12'
$ws.Range('E38').Value = '  +4.04%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.792'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01807'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +1.53%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9006'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.25%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.105'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -2.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9994'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.38'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.45%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '66.96'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.95%  '
$ws.Range('E46').Value = '  -1.65%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.157'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.97%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4018'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.78%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.693'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +1.89%  '
$ws.Range('E50').Value = '  -0.05%  '
$ws.Range('E51').Value = '  +0.23%  '
